{"js": "// Append the new \"Classes\" outline to the end of the document body, after\n// the last existing paragraph (\"Step 9: ... play again at end of game?\").\n//\n// Target shape (see the diff this task was built from):\n//   <w:p/>                                       (blank separator paragraph)\n//   <w:p><w:r><w:t>Classes</w:t></w:r></w:p>\n//   <w:p><w:r><w:t>Game run</w:t></w:r></w:p>\n//   <w:p><w:r><w:t>Player(parent)</w:t></w:r></w:p>\n//   <w:p><w:r><w:t>AI (child player)</w:t></w:r></w:p>\n//   <w:p><w:r><w:t>Human (child player)</w:t></w:r></w:p>\n\nconst body = context.document.body;\n\n// Anchor on the last paragraph currently in the body (\"Step 9: ...\") and\n// insert the new paragraphs after it, one at a time, chaining off the\n// paragraph we just inserted so they land in order, right after the\n// original content and before the (implicit) end of the body.\nlet anchor = body.paragraphs.getLast();\n\nconst newLines = [\n  \"\",                     // blank separator line\n  \"Classes\",\n  \"Game run\",\n  \"Player(parent)\",\n  \"AI (child player)\",\n  \"Human (child player)\"\n];\n\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Append the new \"Classes\" outline to the end of the document body, after\n# the last existing paragraph (\"Step 9: ... play again at end of game?\").\n#\n# Target shape (see the diff this task was built from):\n#   <w:p/>                                       (blank separator paragraph)\n#   <w:p><w:r><w:t>Classes</w:t></w:r></w:p>\n#   <w:p><w:r><w:t>Game run</w:t></w:r></w:p>\n#   <w:p><w:r><w:t>Player(parent)</w:t></w:r></w:p>\n#   <w:p><w:r><w:t>AI (child player)</w:t></w:r></w:p>\n#   <w:p><w:r><w:t>Human (child player)</w:t></w:r></w:p>\n\n$d = $word.ActiveDocument\n\n$lines = @(\"\", \"Classes\", \"Game run\", \"Player(parent)\", \"AI (child player)\", \"Human (child player)\")\n\nforeach ($line in $lines) {\n    # Always re-anchor on the current last paragraph so each new paragraph\n    # lands right after the one before it (and after all original content).\n    $tailRange = $d.Paragraphs.Last.Range\n    $tailRange.InsertParagraphAfter()\n    if ($line -ne \"\") {\n        $newRange = $d.Paragraphs.Last.Range\n        $newRange.Text = $line\n    }\n}\n"}
